$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($r = 3; $r -le 27; $r++) {
    # Column H: "PERIOD TO EXPIRE" - decrement by 1 day
    $hCell = $ws.Cells.Item($r, 8)
    $hCell.Value2 = $hCell.Value2 - 1

    # Column I: "LAST UPDATE" - move from 03-Nov-2025 to 04-Nov-2025.
    # Prefix with an apostrophe so Excel stores it as literal text
    # (matching the source file's text representation) instead of
    # auto-converting the string into a date serial number.
    $iCell = $ws.Cells.Item($r, 9)
    $iCell.Value2 = "'04-Nov-2025"
}
